$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J into the new column K for each data row,
# then fill in the corresponding 2021 values.

# Row 3: bottom border row, no value - formatting only
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 4: header year row
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K4").Value = 2021

# Row 5
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K5").Value = 375

# Row 6
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K6").Value = "-"

# Row 7
$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K7").Value = 5

# Row 8
$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K8").Value = "-"

# Row 9
$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K9").Value = 18

# Row 10
$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K10").Value = 150

$excel.CutCopyMode = 0

# Selection ends on K7, matching the saved view state of the source file.
$ws.Range("K7").Select()
